$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted above the existing row 175,
# shifting every subsequent record (old 175 -> 176, ..., old 200 -> 201)
# down by one row.
$ws.Rows(175).Insert()

# Populate the newly inserted row 175 with the new observation.
$ws.Cells.Item(175, 1).Value = 6
$ws.Cells.Item(175, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(175, 3).Value = "Metropolitana"
$ws.Cells.Item(175, 4).Value = 44522
$ws.Cells.Item(175, 5).Value = 13
$ws.Cells.Item(175, 6).Value = 100112026
$ws.Cells.Item(175, 7).Value = "Haba"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 1000
$ws.Cells.Item(175, 11).Value = 7000
$ws.Cells.Item(175, 12).Value = 8000
$ws.Cells.Item(175, 13).Value = 7440
$ws.Cells.Item(175, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(175, 15).Value = "Región Metropolitana"
$ws.Cells.Item(175, 16).Value = 298
$ws.Cells.Item(175, 17).Value = 25
$ws.Cells.Item(175, 18).Value = "Hortaliza"
